$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 101215
$ws.Range("B6").Value = 91808
$ws.Range("B7").Value = 92267
$ws.Range("B8").Value = 92267
$ws.Range("B9").Value = 91808
$ws.Range("B11").Value = 57830
